# The data rows (2-22) on both worksheets get reshuffled into a new order;
# the contents of each row are unchanged, only their row position moves.
# We use Range.Copy (rather than re-typing .Value) so that cell types
# (e.g. the text-typed UPC columns that look numeric) survive the move
# exactly as they were, instead of Excel re-inferring them as numbers.
#
# Column G (giant_description) is always blank on every row, so it is left
# untouched - no need to move it.

$wb = $excel.ActiveWorkbook

# Maps each NEW row number (2-22) to the OLD row number whose data should
# land there, for both sheets (identical re-ordering on each).
$rowMap = @{
    2  = 3
    3  = 9
    4  = 5
    5  = 16
    6  = 11
    7  = 13
    8  = 6
    9  = 21
    10 = 7
    11 = 22
    12 = 2
    13 = 8
    14 = 17
    15 = 12
    16 = 4
    17 = 18
    18 = 20
    19 = 10
    20 = 14
    21 = 19
    22 = 15
}

foreach ($ws in $wb.Worksheets) {
    # 1) Stage a snapshot of rows 2-22 (columns A:F and H) into a scratch
    #    area well below the used range (rows 1002-1022) so that writes in
    #    step 2 can never clobber data that hasn't been read yet.
    for ($r = 2; $r -le 22; $r++) {
        $stageRow = $r + 1000
        $ws.Range("A$r" + ":F$r").Copy($ws.Range("A$stageRow" + ":F$stageRow"))
        $ws.Range("H$r").Copy($ws.Range("H$stageRow"))
    }

    # 2) Write each row back to its new position from the staged copy.
    foreach ($newRow in $rowMap.Keys) {
        $oldRow = $rowMap[$newRow]
        $stageRow = $oldRow + 1000
        $ws.Range("A$stageRow" + ":F$stageRow").Copy($ws.Range("A$newRow" + ":F$newRow"))
        $ws.Range("H$stageRow").Copy($ws.Range("H$newRow"))
    }

    # 3) Clean up the scratch area.
    $ws.Range("A1002:H1022").Clear()
}
